$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B2:B5 values with the real data
$ws.Range("B2").Value = 24
$ws.Range("B3").Value = 22
$ws.Range("B4").Value = 13
$ws.Range("B5").Value = 16

# Append new rows 6-11 with additional companies and their values
$ws.Range("A6").Value = "e"
$ws.Range("B6").Value = 16

$ws.Range("A7").Value = "f"
$ws.Range("B7").Value = 17

$ws.Range("A8").Value = "g"
$ws.Range("B8").Value = 18

$ws.Range("A9").Value = "h"
$ws.Range("B9").Value = 20

$ws.Range("A10").Value = "i"
$ws.Range("B10").Value = 18

$ws.Range("A11").Value = "average"
$ws.Range("B11").Value = 50

# Move the active selection down to A12, matching post-edit cursor position
$ws.Range("A12").Select()
